$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out all the existing data rows (3-14), we'll rebuild rows 3-20 from scratch
$ws.Range("B3:J14").ClearContents()

# --- Row 3 ---
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = "Num1"
$ws.Range("F3").Value = "Num2"
$ws.Range("G3").Value = "Num3"

# --- Row 4 ---
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = "AVG(Num1, Num2)"
$ws.Range("F4").Value = "Num1"
$ws.Range("G4").Value = "Num2"

# --- Row 5 ---
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = "Num1"
$ws.Range("F5").Value = "Num2"

# --- Row 6 ---
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = "Num1"
$ws.Range("F6").Value = "Num1"
$ws.Range("I6").Value = "Need to review cases like"

# --- Row 7 ---
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = "Num1"
$ws.Range("G7").Value = "Num1"

# --- Row 8 ---
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = "Num1"
$ws.Range("G8").Value = "Num1"

# --- Row 9 ---
$ws.Range("B9").Value = 8
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = "Num1"
$ws.Range("F9").Value = "Num1"

# --- Row 10 ---
$ws.Range("B10").Value = 9
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = "AVG(Num1, Num2)"
$ws.Range("F10").Value = "Num1"
$ws.Range("G10").Value = "Num2"
$ws.Range("I10").Value = "Just 2,000,000/8,000,000 (double comma) which needs sorting. Maybe as a new coding?"

# --- Row 11 ---
$ws.Range("B11").Value = 10
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = "Num1"
$ws.Range("I11").Value = "Case 20, 21 also messing up…"

# --- Row 12 ---
$ws.Range("B12").Value = 11
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = 2
$ws.Range("E12").Value = "AVG(Num1, Num2)"
$ws.Range("F12").Value = "Num1"
$ws.Range("G12").Value = "Num2"

# --- Row 13 ---
$ws.Range("B13").Value = 13
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = 3
$ws.Range("E13").Value = "Num1"
$ws.Range("F13").Value = "Num2"
$ws.Range("G13").Value = "Num3"

# --- Row 14 ---
$ws.Range("B14").Value = 14
$ws.Range("I14").Value = "TODO:"

# --- Row 15 ---
$ws.Range("B15").Value = 15
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = "Num1"
$ws.Range("F15").Value = "Num1"
$ws.Range("G15").Value = "Num1"
$ws.Range("I15").Value = "Handle double commas in case 14 (or as a workaround in case 22, then add)"

# --- Row 16 ---
$ws.Range("B16").Value = 19
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = "AVG(Num1, Num2)"
$ws.Range("F16").Value = "Num1"
$ws.Range("G16").Value = "Num2"
$ws.Range("I16").Value = "Sort out 20, 21 - the regex is crap"

# --- Row 17 ---
$ws.Range("B17").Value = 20
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = "Num1"
$ws.Range("F17").Value = "Num1"
$ws.Range("G17").Value = "Num1"
$ws.Range("I17").Value = "[330 000-370 000]"
$ws.Range("J17").Value = "These cases"

# --- Row 18 ---
$ws.Range("B18").Value = 21
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = "Num1"
$ws.Range("F18").Value = "Num1"
$ws.Range("G18").Value = "Num1"

# --- Row 19 ---
$ws.Range("B19").Value = 22
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = "AVG(Num1, Num2)"
$ws.Range("F19").Value = "Num1"
$ws.Range("G19").Value = "Num2"

# --- Row 20 ---
$ws.Range("B20").Value = 23
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = "AVG(Num1, Num2)"
$ws.Range("F20").Value = "Num1"
$ws.Range("G20").Value = "Num2"

# Update the active selection to match the new state
[void]$ws.Range("E20").Select()
